$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace the text of the "meta description" italic paragraph near the
#    end of the document with the new DALLE image prompt, keeping its
#    run/paragraph formatting (leading empty run + italic run) untouched.
#    NOTE: this must happen before step 3 below duplicates this sentence
#    into a new paragraph near the top of the document.
# ---------------------------------------------------------------------------
$oldDesc = "Explore the Wizard of Oz world with Book of Oz Lock 'n Spin. Get the Lock 'n Spin feature, customizable paylines and play for free."
$descRange = $d.Content
$descFound = $descRange.Find.Execute($oldDesc)
if (-not $descFound) {
    throw "Could not locate the description paragraph to replace"
}

$quoteChar = [char]34
$curlyOpenSingle = [char]0x2018
$newPrompt = "Prompt for DALLE: Create a cartoon-style feature image for " + $quoteChar + "Book of Oz Lock " + $curlyOpenSingle + "N Spins" + $quoteChar + " that features a happy Maya warrior with glasses. The background should be green and the warrior should be holding a magic book with the game's title on it. The warrior should be surrounded by symbols from the game, such as playing card symbols and magic filters in the shape of flowers, hearts, spades, and diamonds. Use bright colors and make the image dynamic and engaging to attract potential players."

$descRange.Text = $newPrompt

# ---------------------------------------------------------------------------
# 2) Remove the duplicated bold "Play Book of Oz Lock 'n Spin Free | Slot
#    Game Review" paragraph that used to sit just above the (now repurposed)
#    italic paragraph near the end of the document.
# ---------------------------------------------------------------------------
$dupRange = $d.Content
$dupRange.Find.Execute("Play Book of Oz Lock 'n Spin Free | Slot Game Review") | Out-Null
$dupRange.Collapse(0)
$dupFound = $dupRange.Find.Execute("Play Book of Oz Lock 'n Spin Free | Slot Game Review")
if (-not $dupFound) {
    throw "Could not locate the duplicated bold title paragraph"
}

$dupPara = $dupRange.Paragraphs(1)
$dupPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Insert a new "Meta description: ..." paragraph right after the
#    Heading 1 title paragraph ("Play Book of Oz Lock 'n Spin Free | Slot
#    Game Review").
# ---------------------------------------------------------------------------
$titleRange = $d.Content
$titleFound = $titleRange.Find.Execute("Play Book of Oz Lock 'n Spin Free | Slot Game Review")
if (-not $titleFound) {
    throw "Could not locate the title paragraph"
}

$insertPos = $titleRange.End
$insertionPoint = $d.Range($insertPos, $insertPos)

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:r/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
    '<w:r><w:t>: Explore the Wizard of Oz world with Book of Oz Lock ' + [char]39 + 'n Spin. Get the Lock ' + [char]39 + 'n Spin feature, customizable paylines and play for free.</w:t></w:r>' +
    '</w:p>'
$insertionPoint.InsertXML($metaXml)
